# Update recomputed precision/recall/F1 statistics (RF and Ensemble rows)
# after re-running the per-class-label classification & evaluation pipeline.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 ("RF")
$ws.Range("B6").Value = 0.8319988464678246
$ws.Range("C6").Value = 0.03516173839978369
$ws.Range("D6").Value = 0.7802938196555218
$ws.Range("E6").Value = 0.8071680919831982
$ws.Range("F6").Value = 0.8511185305540144
$ws.Range("G6").Value = 0.8397558015274496
$ws.Range("H6").Value = 0.8816579886189394
$ws.Range("I6").Value = 0.8129002063521378
$ws.Range("J6").Value = 0.035383199264022
$ws.Range("K6").Value = 0.7772516921453091
$ws.Range("L6").Value = 0.7685225270218438
$ws.Range("M6").Value = 0.8401653260524228
$ws.Range("N6").Value = 0.8182521109516017
$ws.Range("O6").Value = 0.8603093755895114
$ws.Range("P6").Value = 0.7248313897793925
$ws.Range("Q6").Value = 0.02755076777127297
$ws.Range("R6").Value = 0.6956881094979303
$ws.Range("S6").Value = 0.7129774347516282
$ws.Range("T6").Value = 0.7311625505173892
$ws.Range("U6").Value = 0.7092674284998682
$ws.Range("V6").Value = 0.7750614256301468
$ws.Range("W6").Value = 0.7967742207061101
$ws.Range("X6").Value = 0.02527860066961275
$ws.Range("Y6").Value = 0.7546500163990978
$ws.Range("Z6").Value = 0.7887701156904849
$ws.Range("AA6").Value = 0.8161909784203634
$ws.Range("AB6").Value = 0.7961939714061956
$ws.Range("AC6").Value = 0.8280660216144087
$ws.Range("AD6").Value = 0.7949131375952623
$ws.Range("AE6").Value = 0.02277875228516454
$ws.Range("AF6").Value = 0.7778604203815195
$ws.Range("AG6").Value = 0.7867177681387196
$ws.Range("AH6").Value = 0.799009004760278
$ws.Range("AI6").Value = 0.7738987882612498
$ws.Range("AJ6").Value = 0.837079706434545
$ws.Range("AK6").Value = 0.8213291381248924
$ws.Range("AL6").Value = 0.04885307290215209
$ws.Range("AM6").Value = 0.7849891966003388
$ws.Range("AN6").Value = 0.7522637176745373
$ws.Range("AO6").Value = 0.8501242616112799
$ws.Range("AP6").Value = 0.8273947897037879
$ws.Range("AQ6").Value = 0.8918737250345183

# Row 7 ("Ensemble")
$ws.Range("B7").Value = 0.8438845433216151
$ws.Range("C7").Value = 0.04876288202835144
$ws.Range("D7").Value = 0.8410545095875442
$ws.Range("F7").Value = 0.8701572507599673
$ws.Range("I7").Value = 0.8573795959786967
$ws.Range("J7").Value = 0.03144225869023921
$ws.Range("M7").Value = 0.8585770860399209
$ws.Range("N7").Value = 0.8598535593028118
$ws.Range("P7").Value = 0.8293527408589704
$ws.Range("Q7").Value = 0.04604049327081808
$ws.Range("R7").Value = 0.8116584982766393
$ws.Range("S7").Value = 0.7870787846806523
$ws.Range("T7").Value = 0.8398903417658107
$ws.Range("V7").Value = 0.9139250714388401
$ws.Range("W7").Value = 0.8455985765736322
$ws.Range("X7").Value = 0.03318404238425834
$ws.Range("Y7").Value = 0.8411041475390205
$ws.Range("AB7").Value = 0.8482820220593374
$ws.Range("AC7").Value = 0.8699613636049927
$ws.Range("AD7").Value = 0.856295760917458
$ws.Range("AE7").Value = 0.03516062000954683
$ws.Range("AF7").Value = 0.862432265743865
$ws.Range("AG7").Value = 0.7960737244268654
$ws.Range("AH7").Value = 0.8713783648382389
$ws.Range("AJ7").Value = 0.9034458820749144
$ws.Range("AK7").Value = 0.837194762515152
$ws.Range("AL7").Value = 0.02491522337540308
$ws.Range("AM7").Value = 0.8311148269114019
$ws.Range("AN7").Value = 0.8062734756283143
$ws.Range("AQ7").Value = 0.8817315306820649
